# Add a "Result" column (M) to the BillingInfoValidation sheet, marking
# PASS for the rows whose RunMode is "Yes" (rows 2, 4 and 6), matching the
# 29-April discussion about constants and writing results back to Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled the same as the other header cells (bold).
$ws.Range("M1").Value = "Result"
$ws.Range("M1").Font.Bold = $true

# Result values for the rows that were validated.
$ws.Range("M2").Value = "PASS"
$ws.Range("M4").Value = "PASS"
$ws.Range("M6").Value = "PASS"

# Leave the active selection on the newly added header cell.
$ws.Range("M1").Select() | Out-Null
